$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 960.2692
$ws.Range("I15").Value = 960.2692
$ws.Range("K15").Value = 2880.8076
$ws.Range("M15").Value = -2711.8076
$ws.Range("H86").Value = 10000
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 10000
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H131").Value = 4999
$ws.Range("I131").Value = 4998
$ws.Range("K131").Value = 14994
$ws.Range("M131").Value = -9954
$ws.Range("H132").Value = 2315.2666
$ws.Range("I132").Value = 1030
$ws.Range("K132").Value = 3090
$ws.Range("M132").Value = -560
$ws.Range("H137").Value = 3181.7058
$ws.Range("I137").Value = 2624.8333
$ws.Range("J137").Value = 3485.4546
$ws.Range("K137").Value = 7874.499899999999
$ws.Range("L137").Value = 10456.3638
$ws.Range("M137").Value = -5324.499899999999
$ws.Range("N137").Value = -15556.3638

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 500
$ws.Range("I5").Value = 500
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 500
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -388
$ws.Range("N5").ClearContents()
$ws.Range("H32").Value = 7306.1
$ws.Range("I32").Value = 6480.1055
$ws.Range("K32").Value = 6480.1055
$ws.Range("M32").Value = -6193.1055
$ws.Range("H61").Value = 3257
$ws.Range("I61").Value = 3164.125
$ws.Range("K61").Value = 3164.125
$ws.Range("M61").Value = -2952.125
$ws.Range("H63").Value = 10745
$ws.Range("J63").Value = 10745
$ws.Range("L63").Value = 10745
$ws.Range("N63").Value = -12117
$ws.Range("H66").Value = 10745
$ws.Range("J66").Value = 10745
$ws.Range("L66").Value = 53725
$ws.Range("N66").Value = -60589

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -385
$ws.Range("N4").ClearContents()
$ws.Range("H105").Value = 11143.714
$ws.Range("I105").Value = 21603.334
$ws.Range("J105").Value = 3299
$ws.Range("K105").Value = 21603.334
$ws.Range("L105").Value = 3299
$ws.Range("M105").Value = -19856.334
$ws.Range("N105").Value = -6793

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I3").Value = 9999
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 9999
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -9886
$ws.Range("N3").ClearContents()
$ws.Range("H31").Value = 2781.1428
$ws.Range("I31").Value = 2240
$ws.Range("K31").Value = 2240
$ws.Range("M31").Value = -1945
$ws.Range("H34").Value = 2781.1428
$ws.Range("I34").Value = 2240
$ws.Range("K34").Value = 2240
$ws.Range("M34").Value = -2038
$ws.Range("H105").Value = 2356
$ws.Range("I105").Value = 2415.6667
$ws.Range("J105").Value = 1998
$ws.Range("K105").Value = 2415.6667
$ws.Range("L105").Value = 1998
$ws.Range("M105").Value = -668.6667000000002
$ws.Range("N105").Value = -5492

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1999.5
$ws.Range("J68").Value = 1999.5
$ws.Range("L68").Value = 5998.5
$ws.Range("N68").Value = -7620.5
$ws.Range("H71").Value = 1999.5
$ws.Range("J71").Value = 1999.5
$ws.Range("L71").Value = 17995.5
$ws.Range("N71").Value = -26107.5
$ws.Range("H99").Value = 4640.44
$ws.Range("J99").Value = 4828.75
$ws.Range("L99").Value = 14486.25
$ws.Range("N99").Value = -18978.25
$ws.Range("H103").Value = 2884.6667
$ws.Range("I103").Value = 2032
$ws.Range("K103").Value = 6096
$ws.Range("M103").Value = -5217
$ws.Range("H113").Value = 2500.6
$ws.Range("J113").Value = 2875.75
$ws.Range("L113").Value = 8627.25
$ws.Range("N113").Value = -12967.25
$ws.Range("H131").Value = 1644.1818
$ws.Range("J131").Value = 1649.5
$ws.Range("L131").Value = 4948.5
$ws.Range("N131").Value = -15028.5
$ws.Range("H132").Value = 9282.429
$ws.Range("J132").Value = 11395.4
$ws.Range("L132").Value = 102558.6
$ws.Range("N132").Value = -107618.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 772
$ws.Range("I102").Value = 772
$ws.Range("K102").Value = 772
$ws.Range("M102").Value = 850
$ws.Range("H113").Value = 2387.1428
$ws.Range("I113").Value = 2677.5
$ws.Range("K113").Value = 2677.5
$ws.Range("M113").Value = -507.5
$ws.Range("H126").Value = 3049
$ws.Range("I126").Value = 2891
$ws.Range("K126").Value = 8673
$ws.Range("M126").Value = -6203
$ws.Range("H132").Value = 3995.9375
$ws.Range("I132").Value = 3494.3635
$ws.Range("J132").Value = 5099.4
$ws.Range("K132").Value = 10483.0905
$ws.Range("L132").Value = 15298.2
$ws.Range("M132").Value = -7953.0905
$ws.Range("N132").Value = -20358.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2224.8333
$ws.Range("I40").Value = 2224.8333
$ws.Range("K40").Value = 2224.8333
$ws.Range("M40").Value = -2088.8333
$ws.Range("H46").Value = 3487.5
$ws.Range("I46").Value = 3500
$ws.Range("J46").Value = 3483.3333
$ws.Range("K46").Value = 3500
$ws.Range("L46").Value = 3483.3333
$ws.Range("M46").Value = -3312
$ws.Range("N46").Value = -3859.3333
$ws.Range("H55").Value = 512
$ws.Range("I55").Value = 597
$ws.Range("J55").Value = 235.75
$ws.Range("K55").Value = 597
$ws.Range("L55").Value = 235.75
$ws.Range("M55").Value = -424
$ws.Range("N55").Value = -581.75
$ws.Range("H82").Value = 4133.25
$ws.Range("I82").Value = 1638
$ws.Range("J82").Value = 9123.75
$ws.Range("K82").Value = 1638
$ws.Range("L82").Value = 9123.75
$ws.Range("M82").Value = -1277
$ws.Range("N82").Value = -9845.75
$ws.Range("H85").Value = 4133.25
$ws.Range("I85").Value = 1638
$ws.Range("J85").Value = 9123.75
$ws.Range("K85").Value = 1638
$ws.Range("L85").Value = 9123.75
$ws.Range("M85").Value = -390
$ws.Range("N85").Value = -11619.75
$ws.Range("H132").Value = 5128.143
$ws.Range("I132").Value = 3966.3333
$ws.Range("K132").Value = 11898.9999
$ws.Range("M132").Value = -9368.999899999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 55000
$ws.Range("J109").Value = 55000
$ws.Range("L109").Value = 55000
$ws.Range("N109").Value = -57774

